$wb = $excel.ActiveWorkbook

# --- Sheet: y_fitted_on_begin_2016 ---
$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws1.Cells.Item(2,1).Value = 1981
$ws1.Cells.Item(2,2).Value = 16.69374857514662
$ws1.Cells.Item(3,1).Value = 1982
$ws1.Cells.Item(3,2).Value = 16.69227074751712
$ws1.Cells.Item(4,1).Value = 1983
$ws1.Cells.Item(4,2).Value = 16.41874885370408
$ws1.Cells.Item(5,1).Value = 1984
$ws1.Cells.Item(5,2).Value = 16.6039372162313
$ws1.Cells.Item(6,1).Value = 1985
$ws1.Cells.Item(6,2).Value = 16.60876429117652
$ws1.Cells.Item(7,1).Value = 1986
$ws1.Cells.Item(7,2).Value = 16.42611259825074
$ws1.Cells.Item(8,1).Value = 1987
$ws1.Cells.Item(8,2).Value = 16.27707519156587
$ws1.Cells.Item(9,1).Value = 1988
$ws1.Cells.Item(9,2).Value = 16.23800579942649
$ws1.Cells.Item(10,1).Value = 1989
$ws1.Cells.Item(10,2).Value = 16.16313926321722
$ws1.Cells.Item(11,1).Value = 1990
$ws1.Cells.Item(11,2).Value = 16.69886245704583
$ws1.Cells.Item(12,1).Value = 1991
$ws1.Cells.Item(12,2).Value = 16.49485117600963
$ws1.Cells.Item(13,1).Value = 1992
$ws1.Cells.Item(13,2).Value = 15.91395673294453
$ws1.Cells.Item(14,1).Value = 1993
$ws1.Cells.Item(14,2).Value = 16.01605675374332
$ws1.Cells.Item(15,1).Value = 1994
$ws1.Cells.Item(15,2).Value = 16.10968983573308
$ws1.Cells.Item(16,1).Value = 1995
$ws1.Cells.Item(16,2).Value = 16.63298357721069
$ws1.Cells.Item(17,1).Value = 1996
$ws1.Cells.Item(17,2).Value = 17.35763124745286
$ws1.Cells.Item(18,1).Value = 1997
$ws1.Cells.Item(18,2).Value = 17.19828698199613
$ws1.Cells.Item(19,1).Value = 1998
$ws1.Cells.Item(19,2).Value = 16.47294399445137
$ws1.Cells.Item(20,1).Value = 1999
$ws1.Cells.Item(20,2).Value = 15.92602433398816
$ws1.Cells.Item(21,1).Value = 2000
$ws1.Cells.Item(21,2).Value = 15.55858982949081
$ws1.Cells.Item(22,1).Value = 2001
$ws1.Cells.Item(22,2).Value = 16.06270268919287
$ws1.Cells.Item(23,1).Value = 2002
$ws1.Cells.Item(23,2).Value = 15.64938837375653
$ws1.Cells.Item(24,1).Value = 2003
$ws1.Cells.Item(24,2).Value = 15.81246671767336
$ws1.Cells.Item(25,1).Value = 2004
$ws1.Cells.Item(25,2).Value = 15.82277580339271
$ws1.Cells.Item(26,1).Value = 2005
$ws1.Cells.Item(26,2).Value = 15.95006300571342
$ws1.Cells.Item(27,1).Value = 2006
$ws1.Cells.Item(27,2).Value = 16.03114521591096
$ws1.Cells.Item(28,1).Value = 2007
$ws1.Cells.Item(28,2).Value = 16.91918165897853
$ws1.Cells.Item(29,1).Value = 2008
$ws1.Cells.Item(29,2).Value = 16.58450706259895
$ws1.Cells.Item(30,1).Value = 2009
$ws1.Cells.Item(30,2).Value = 16.66877061872903
$ws1.Cells.Item(31,1).Value = 2010
$ws1.Cells.Item(31,2).Value = 16.69928793282159
$ws1.Cells.Item(32,1).Value = 2011
$ws1.Cells.Item(32,2).Value = 16.50961644307383
$ws1.Cells.Item(33,1).Value = 2012
$ws1.Cells.Item(33,2).Value = 15.88267280339234
$ws1.Cells.Item(34,1).Value = 2013
$ws1.Cells.Item(34,2).Value = 15.69853876525039
$ws1.Cells.Item(35,1).Value = 2014
$ws1.Cells.Item(35,2).Value = 15.27080496949527
$ws1.Cells.Item(36,1).Value = 2015
$ws1.Cells.Item(36,2).Value = 15.1102468738348
$ws1.Cells.Item(37,1).Value = 2016
$ws1.Cells.Item(37,2).Value = 15.46925328829514

# --- Sheet: y_pred_on_2017_2021 ---
$ws2 = $wb.Worksheets.Item("y_pred_on_2017_2021")
$ws2.Cells.Item(2,2).Value = 15.16538317643262
$ws2.Cells.Item(3,2).Value = 15.15073183569228
$ws2.Cells.Item(4,2).Value = 15.12395466891853
$ws2.Cells.Item(5,2).Value = 15.08898887598999
$ws2.Cells.Item(6,2).Value = 15.04857248477361

# --- Sheet: y_fitted_on_begin_2021 ---
$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws3.Cells.Item(2,1).Value = 1981
$ws3.Cells.Item(2,2).Value = 16.68457877835144
$ws3.Cells.Item(3,1).Value = 1982
$ws3.Cells.Item(3,2).Value = 16.68419575185602
$ws3.Cells.Item(4,1).Value = 1983
$ws3.Cells.Item(4,2).Value = 16.39467212405009
$ws3.Cells.Item(5,1).Value = 1984
$ws3.Cells.Item(5,2).Value = 16.59588602824662
$ws3.Cells.Item(6,1).Value = 1985
$ws3.Cells.Item(6,2).Value = 16.59757862640502
$ws3.Cells.Item(7,1).Value = 1986
$ws3.Cells.Item(7,2).Value = 16.41328604270906
$ws3.Cells.Item(8,1).Value = 1987
$ws3.Cells.Item(8,2).Value = 16.2523562599051
$ws3.Cells.Item(9,1).Value = 1988
$ws3.Cells.Item(9,2).Value = 16.20378762390939
$ws3.Cells.Item(10,1).Value = 1989
$ws3.Cells.Item(10,2).Value = 16.1293165396719
$ws3.Cells.Item(11,1).Value = 1990
$ws3.Cells.Item(11,2).Value = 16.71271715232487
$ws3.Cells.Item(12,1).Value = 1991
$ws3.Cells.Item(12,2).Value = 16.50055693646327
$ws3.Cells.Item(13,1).Value = 1992
$ws3.Cells.Item(13,2).Value = 15.89860013263793
$ws3.Cells.Item(14,1).Value = 1993
$ws3.Cells.Item(14,2).Value = 16.02400038161839
$ws3.Cells.Item(15,1).Value = 1994
$ws3.Cells.Item(15,2).Value = 16.13939277658897
$ws3.Cells.Item(16,1).Value = 1995
$ws3.Cells.Item(16,2).Value = 16.68991978232189
$ws3.Cells.Item(17,1).Value = 1996
$ws3.Cells.Item(17,2).Value = 17.4712210737991
$ws3.Cells.Item(18,1).Value = 1997
$ws3.Cells.Item(18,2).Value = 17.31568837922985
$ws3.Cells.Item(19,1).Value = 1998
$ws3.Cells.Item(19,2).Value = 16.54130919410686
$ws3.Cells.Item(20,1).Value = 1999
$ws3.Cells.Item(20,2).Value = 15.95558275149642
$ws3.Cells.Item(21,1).Value = 2000
$ws3.Cells.Item(21,2).Value = 15.5579915886593
$ws3.Cells.Item(22,1).Value = 2001
$ws3.Cells.Item(22,2).Value = 16.10402009686713
$ws3.Cells.Item(23,1).Value = 2002
$ws3.Cells.Item(23,2).Value = 15.6579306901002
$ws3.Cells.Item(24,1).Value = 2003
$ws3.Cells.Item(24,2).Value = 15.82889537621459
$ws3.Cells.Item(25,1).Value = 2004
$ws3.Cells.Item(25,2).Value = 15.83391987153613
$ws3.Cells.Item(26,1).Value = 2005
$ws3.Cells.Item(26,2).Value = 15.95400783524422
$ws3.Cells.Item(27,1).Value = 2006
$ws3.Cells.Item(27,2).Value = 16.03148803871925
$ws3.Cells.Item(28,1).Value = 2007
$ws3.Cells.Item(28,2).Value = 16.98423779974829
$ws3.Cells.Item(29,1).Value = 2008
$ws3.Cells.Item(29,2).Value = 16.60406325767243
$ws3.Cells.Item(30,1).Value = 2009
$ws3.Cells.Item(30,2).Value = 16.6606372271191
$ws3.Cells.Item(31,1).Value = 2010
$ws3.Cells.Item(31,2).Value = 16.69170857561161
$ws3.Cells.Item(32,1).Value = 2011
$ws3.Cells.Item(32,2).Value = 16.43554308289066
$ws3.Cells.Item(33,1).Value = 2012
$ws3.Cells.Item(33,2).Value = 15.741871537443
$ws3.Cells.Item(34,1).Value = 2013
$ws3.Cells.Item(34,2).Value = 15.55267746232432
$ws3.Cells.Item(35,1).Value = 2014
$ws3.Cells.Item(35,2).Value = 15.06501181894495
$ws3.Cells.Item(36,1).Value = 2015
$ws3.Cells.Item(36,2).Value = 14.87157744692874
$ws3.Cells.Item(37,1).Value = 2016
$ws3.Cells.Item(37,2).Value = 15.25802119504917
$ws3.Cells.Item(38,1).Value = 2017
$ws3.Cells.Item(38,2).Value = 14.91621529978475
$ws3.Cells.Item(39,1).Value = 2018
$ws3.Cells.Item(39,2).Value = 14.58399189597616
$ws3.Cells.Item(40,1).Value = 2019
$ws3.Cells.Item(40,2).Value = 14.35049850627647
$ws3.Cells.Item(41,1).Value = 2020
$ws3.Cells.Item(41,2).Value = 13.63583590498246
$ws3.Cells.Item(42,1).Value = 2021
$ws3.Cells.Item(42,2).Value = 14.24530912303092
$ws3.Rows.Item(43).Delete() | Out-Null

# --- Sheet: y_pred_on_2022_2026 ---
$ws4 = $wb.Worksheets.Item("y_pred_on_2022_2026")
$ws4.Cells.Item(2,2).Value = 14.98554761906389
$ws4.Cells.Item(3,2).Value = 15.84730033002362
$ws4.Cells.Item(4,2).Value = 16.60219084422058
$ws4.Cells.Item(5,2).Value = 17.21096277148818
$ws4.Cells.Item(6,2).Value = 17.64944078431874

Write-Output "done"